$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.004809333333333
$ws.Range("H2").Value = 12.014428
$ws.Range("I2").Value = 0.04337108182100274
$ws.Range("J2").Value = 0.04337108182100273
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 163.7119853333333
$ws.Range("N2").Value = 491.135956
$ws.Range("O2").Value = 0.2754003062401033
$ws.Range("P2").Value = 0.2754003062401033
$ws.Range("Q2").Value = 655.6352868414631
$ws.Range("R2").Value = 5900.717581573167
$ws.Range("S2").Value = 0.01194440921546873
$ws.Range("T2").Value = 0.01194440921546873

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.004809333333333
$ws.Range("H3").Value = 12.014428
$ws.Range("I3").Value = 0.04337108182100274
$ws.Range("J3").Value = 0.04337108182100273
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 172.558497
$ws.Range("N3").Value = 517.675491
$ws.Range("O3").Value = 0.290282124557779
$ws.Range("P3").Value = 0.290282124557779
$ws.Range("Q3").Value = 691.0638793315719
$ws.Range("R3").Value = 6219.574913984147
$ws.Range("S3").Value = 0.01258984977536994
$ws.Range("T3").Value = 0.01258984977536994

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.004809333333333
$ws.Range("H4").Value = 12.014428
$ws.Range("I4").Value = 0.04337108182100274
$ws.Range("J4").Value = 0.04337108182100273
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 111.4881643333333
$ws.Range("N4").Value = 334.464493
$ws.Range("O4").Value = 0.1875481171218523
$ws.Range("P4").Value = 0.1875481171218523
$ws.Range("Q4").Value = 446.4888410783337
$ws.Range("R4").Value = 4018.399569705004
$ws.Range("S4").Value = 0.00813416473306686
$ws.Range("T4").Value = 0.00813416473306686

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.004809333333333
$ws.Range("H5").Value = 12.014428
$ws.Range("I5").Value = 0.04337108182100274
$ws.Range("J5").Value = 0.04337108182100273
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 86.95798233333333
$ws.Range("N5").Value = 260.873947
$ws.Range("O5").Value = 0.1462828449356383
$ws.Range("P5").Value = 0.1462828449356383
$ws.Range("Q5").Value = 348.2501392563684
$ws.Range("R5").Value = 3134.251253307315
$ws.Range("S5").Value = 0.006344445236712626
$ws.Range("T5").Value = 0.006344445236712625

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4.004809333333333
$ws.Range("H6").Value = 12.014428
$ws.Range("I6").Value = 0.04337108182100274
$ws.Range("J6").Value = 0.04337108182100273
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 59.73436333333333
$ws.Range("N6").Value = 179.20309
$ws.Range("O6").Value = 0.100486607144627
$ws.Range("P6").Value = 0.100486607144627
$ws.Range("Q6").Value = 239.2247357980577
$ws.Range("R6").Value = 2153.02262218252
$ws.Range("S6").Value = 0.004358212860384575
$ws.Range("T6").Value = 0.004358212860384574

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 64.80903633333334
$ws.Range("H7").Value = 194.427109
$ws.Range("I7").Value = 0.7018656279483316
$ws.Range("J7").Value = 0.7018656279483316
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 163.7119853333333
$ws.Range("N7").Value = 491.135956
$ws.Range("O7").Value = 0.2754003062401033
$ws.Range("P7").Value = 0.2754003062401033
$ws.Range("Q7").Value = 10610.01600567014
$ws.Range("R7").Value = 95490.1440510312
$ws.Range("S7").Value = 0.1932940088763729
$ws.Range("T7").Value = 0.1932940088763729

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 64.80903633333334
$ws.Range("H8").Value = 194.427109
$ws.Range("I8").Value = 0.7018656279483316
$ws.Range("J8").Value = 0.7018656279483316
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 172.558497
$ws.Range("N8").Value = 517.675491
$ws.Range("O8").Value = 0.290282124557779
$ws.Range("P8").Value = 0.290282124557779
$ws.Range("Q8").Value = 11183.34990169839
$ws.Range("R8").Value = 100650.1491152855
$ws.Range("S8").Value = 0.2037390456349213
$ws.Range("T8").Value = 0.2037390456349213

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 64.80903633333334
$ws.Range("H9").Value = 194.427109
$ws.Range("I9").Value = 0.7018656279483316
$ws.Range("J9").Value = 0.7018656279483316
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 111.4881643333333
$ws.Range("N9").Value = 334.464493
$ws.Range("O9").Value = 0.1875481171218523
$ws.Range("P9").Value = 0.1875481171218523
$ws.Range("Q9").Value = 7225.440493015638
$ws.Range("R9").Value = 65028.96443714074
$ws.Range("S9").Value = 0.1316335769942561
$ws.Range("T9").Value = 0.1316335769942561

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 64.80903633333334
$ws.Range("H10").Value = 194.427109
$ws.Range("I10").Value = 0.7018656279483316
$ws.Range("J10").Value = 0.7018656279483316
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 86.95798233333333
$ws.Range("N10").Value = 260.873947
$ws.Range("O10").Value = 0.1462828449356383
$ws.Range("P10").Value = 0.1462828449356383
$ws.Range("Q10").Value = 5635.663036514358
$ws.Range("R10").Value = 50720.96732862922
$ws.Range("S10").Value = 0.1026709008188202
$ws.Range("T10").Value = 0.1026709008188202

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 64.80903633333334
$ws.Range("H11").Value = 194.427109
$ws.Range("I11").Value = 0.7018656279483316
$ws.Range("J11").Value = 0.7018656279483316
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 59.73436333333333
$ws.Range("N11").Value = 179.20309
$ws.Range("O11").Value = 0.100486607144627
$ws.Range("P11").Value = 0.100486607144627
$ws.Range("Q11").Value = 3871.326523618534
$ws.Range("R11").Value = 34841.93871256681
$ws.Range("S11").Value = 0.07052809562396091
$ws.Range("T11").Value = 0.07052809562396091

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 5.472599666666667
$ws.Range("H12").Value = 16.417799
$ws.Range("I12").Value = 0.05926688342963785
$ws.Range("J12").Value = 0.05926688342963785
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 163.7119853333333
$ws.Range("N12").Value = 491.135956
$ws.Range("O12").Value = 0.2754003062401033
$ws.Range("P12").Value = 0.2754003062401033
$ws.Range("Q12").Value = 895.9301563645383
$ws.Range("R12").Value = 8063.371407280844
$ws.Range("S12").Value = 0.01632211784641877
$ws.Range("T12").Value = 0.01632211784641877

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 5.472599666666667
$ws.Range("H13").Value = 16.417799
$ws.Range("I13").Value = 0.05926688342963785
$ws.Range("J13").Value = 0.05926688342963785
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 172.558497
$ws.Range("N13").Value = 517.675491
$ws.Range("O13").Value = 0.290282124557779
$ws.Range("P13").Value = 0.290282124557779
$ws.Range("Q13").Value = 944.3435731627011
$ws.Range("R13").Value = 8499.09215846431
$ws.Range("S13").Value = 0.0172041168378735
$ws.Range("T13").Value = 0.0172041168378735

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 5.472599666666667
$ws.Range("H14").Value = 16.417799
$ws.Range("I14").Value = 0.05926688342963785
$ws.Range("J14").Value = 0.05926688342963785
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 111.4881643333333
$ws.Range("N14").Value = 334.464493
$ws.Range("O14").Value = 0.1875481171218523
$ws.Range("P14").Value = 0.1875481171218523
$ws.Range("Q14").Value = 610.1300909678786
$ws.Range("R14").Value = 5491.170818710908
$ws.Range("S14").Value = 0.01111539239490889
$ws.Range("T14").Value = 0.01111539239490889

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 5.472599666666667
$ws.Range("H15").Value = 16.417799
$ws.Range("I15").Value = 0.05926688342963785
$ws.Range("J15").Value = 0.05926688342963785
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 86.95798233333333
$ws.Range("N15").Value = 260.873947
$ws.Range("O15").Value = 0.1462828449356383
$ws.Range("P15").Value = 0.1462828449356383
$ws.Range("Q15").Value = 475.886225131406
$ws.Range("R15").Value = 4282.976026182653
$ws.Range("S15").Value = 0.008669728318556268
$ws.Range("T15").Value = 0.008669728318556266

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 5.472599666666667
$ws.Range("H16").Value = 16.417799
$ws.Range("I16").Value = 0.05926688342963785
$ws.Range("J16").Value = 0.05926688342963785
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 59.73436333333333
$ws.Range("N16").Value = 179.20309
$ws.Range("O16").Value = 0.100486607144627
$ws.Range("P16").Value = 0.100486607144627
$ws.Range("Q16").Value = 326.9022568665455
$ws.Range("R16").Value = 2942.12031179891
$ws.Range("S16").Value = 0.005955528031880421
$ws.Range("T16").Value = 0.00595552803188042

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 4.689392000000001
$ws.Range("H17").Value = 14.068176
$ws.Range("I17").Value = 0.05078494060376966
$ws.Range("J17").Value = 0.05078494060376965
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 163.7119853333333
$ws.Range("N17").Value = 491.135956
$ws.Range("O17").Value = 0.2754003062401033
$ws.Range("P17").Value = 0.2754003062401033
$ws.Range("Q17").Value = 767.7096743262508
$ws.Range("R17").Value = 6909.387068936256
$ws.Range("S17").Value = 0.01398618819466362
$ws.Range("T17").Value = 0.01398618819466362

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 4.689392000000001
$ws.Range("H18").Value = 14.068176
$ws.Range("I18").Value = 0.05078494060376966
$ws.Range("J18").Value = 0.05078494060376965
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 172.558497
$ws.Range("N18").Value = 517.675491
$ws.Range("O18").Value = 0.290282124557779
$ws.Range("P18").Value = 0.290282124557779
$ws.Range("Q18").Value = 809.194435363824
$ws.Range("R18").Value = 7282.749918274416
$ws.Range("S18").Value = 0.01474196045400287
$ws.Range("T18").Value = 0.01474196045400287

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 4.689392000000001
$ws.Range("H19").Value = 14.068176
$ws.Range("I19").Value = 0.05078494060376966
$ws.Range("J19").Value = 0.05078494060376965
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 111.4881643333333
$ws.Range("N19").Value = 334.464493
$ws.Range("O19").Value = 0.1875481171218523
$ws.Range("P19").Value = 0.1875481171218523
$ws.Range("Q19").Value = 522.8117059194187
$ws.Range("R19").Value = 4705.305353274768
$ws.Range("S19").Value = 0.009524619988382103
$ws.Range("T19").Value = 0.009524619988382105

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 4.689392000000001
$ws.Range("H20").Value = 14.068176
$ws.Range("I20").Value = 0.05078494060376966
$ws.Range("J20").Value = 0.05078494060376965
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 86.95798233333333
$ws.Range("N20").Value = 260.873947
$ws.Range("O20").Value = 0.1462828449356383
$ws.Range("P20").Value = 0.1462828449356383
$ws.Range("Q20").Value = 407.7800666900747
$ws.Range("R20").Value = 3670.020600210672
$ws.Range("S20").Value = 0.00742896559140684
$ws.Range("T20").Value = 0.007428965591406839

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 4.689392000000001
$ws.Range("H21").Value = 14.068176
$ws.Range("I21").Value = 0.05078494060376966
$ws.Range("J21").Value = 0.05078494060376965
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 59.73436333333333
$ws.Range("N21").Value = 179.20309
$ws.Range("O21").Value = 0.100486607144627
$ws.Range("P21").Value = 0.100486607144627
$ws.Range("Q21").Value = 280.1178455404267
$ws.Range("R21").Value = 2521.06060986384
$ws.Range("S21").Value = 0.005103206375314216
$ws.Range("T21").Value = 0.005103206375314216

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 13.362402
$ws.Range("H22").Value = 40.087206
$ws.Range("I22").Value = 0.1447114661972582
$ws.Range("J22").Value = 0.1447114661972582
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 163.7119853333333
$ws.Range("N22").Value = 491.135956
$ws.Range("O22").Value = 0.2754003062401033
$ws.Range("P22").Value = 0.2754003062401033
$ws.Range("Q22").Value = 2187.585360242104
$ws.Range("R22").Value = 19688.26824217894
$ws.Range("S22").Value = 0.03985358210717926
$ws.Range("T22").Value = 0.03985358210717926

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 13.362402
$ws.Range("H23").Value = 40.087206
$ws.Range("I23").Value = 0.1447114661972582
$ws.Range("J23").Value = 0.1447114661972582
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 172.558497
$ws.Range("N23").Value = 517.675491
$ws.Range("O23").Value = 0.290282124557779
$ws.Range("P23").Value = 0.290282124557779
$ws.Range("Q23").Value = 2305.796005429794
$ws.Range("R23").Value = 20752.16404886814
$ws.Range("S23").Value = 0.04200715185561132
$ws.Range("T23").Value = 0.04200715185561132

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 13.362402
$ws.Range("H24").Value = 40.087206
$ws.Range("I24").Value = 0.1447114661972582
$ws.Range("J24").Value = 0.1447114661972582
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 111.4881643333333
$ws.Range("N24").Value = 334.464493
$ws.Range("O24").Value = 0.1875481171218523
$ws.Range("P24").Value = 0.1875481171218523
$ws.Range("Q24").Value = 1489.749670064062
$ws.Range("R24").Value = 13407.74703057656
$ws.Range("S24").Value = 0.02714036301123834
$ws.Range("T24").Value = 0.02714036301123835

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 13.362402
$ws.Range("H25").Value = 40.087206
$ws.Range("I25").Value = 0.1447114661972582
$ws.Range("J25").Value = 0.1447114661972582
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 86.95798233333333
$ws.Range("N25").Value = 260.873947
$ws.Range("O25").Value = 0.1462828449356383
$ws.Range("P25").Value = 0.1462828449356383
$ws.Range("Q25").Value = 1161.967517046898
$ws.Range("R25").Value = 10457.70765342208
$ws.Range("S25").Value = 0.02116880497014239
$ws.Range("T25").Value = 0.02116880497014239

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 13.362402
$ws.Range("H26").Value = 40.087206
$ws.Range("I26").Value = 0.1447114661972582
$ws.Range("J26").Value = 0.1447114661972582
$ws.Range("K26").Value = 3
$ws.Range("M26").Value = 59.73436333333333
$ws.Range("N26").Value = 179.20309
$ws.Range("O26").Value = 0.100486607144627
$ws.Range("P26").Value = 0.100486607144627
$ws.Range("Q26").Value = 798.19457607406
$ws.Range("R26").Value = 7183.751184666539
$ws.Range("S26").Value = 0.01454156425308685
$ws.Range("T26").Value = 0.01454156425308685
